$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Append a "?" to the end of the solution-specification question (cell C34).
$ws.Range("C34").Value = "A domain-specific diagram or form for explanation of the solution specification should be provided. Right now we just use a behavior definition document. Is that good enough?"

# 2. New update rows 42-45 documenting the new problem-diagram / sub-argument work.
$ws.Range("G42").Value = 45478
$ws.Range("H42").Value = "Rebuild problem space diagram to be more accurate."
$ws.Range("I42").Value = "OPEN"
$ws.Rows.Item(42).RowHeight = 17

$ws.Range("G43").Value = 45478
$ws.Range("H43").Value = "Applied a state-based model to the definition of each problem and showed how they are defined against requirements in various hierarchies of sub-solution."
$ws.Range("I43").Value = "OPEN"
$ws.Rows.Item(43).RowHeight = 34

$ws.Range("G44").Value = 45478
$ws.Range("H44").Value = "A version of the argument that refers to the problem diagram and explains intermediate solutions and the identification of the specific lower-level problem to solve has been written"
$ws.Range("I44").Value = "OPEN"
$ws.Rows.Item(44).RowHeight = 51

$ws.Range("H45").Value = "Discussion and comparison of the two problem arguent versions (one that refers to problemSpace.png, and another that cites it and explains the diagram reasoning.)"
$ws.Rows.Item(45).RowHeight = 34

# 3. Issue 9's status (B33) changed from OPEN to CANCELLED, superseded by the new problem-diagram work.
$ws.Range("B33").Value = "CANCELLED"

# 4. Update the selection to match the saved cursor position.
$ws.Range("C34").Select()
